$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.264.86'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.620.05'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.93'
$ws.Range('E6').Value = '  -1.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3787'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '52.03'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3598'
$ws.Range('E9').Value = '  -2.38%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08054'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('B11').Value = 'BinanceUSD'
$ws.Range('C11').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.003'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.218'
$ws.Range('E12').Value = '  -5.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.50'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.523'
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('E15').Value = '  -4.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.195'
$ws.Range('E16').Value = '  -3.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.622.35'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.35'
$ws.Range('E18').Value = '  -1.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06902'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.85'
$ws.Range('E20').Value = '  -3.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.401'
$ws.Range('E22').Value = '  -3.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '23.277.06'
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.67'
$ws.Range('E24').Value = '  -2.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.162'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.455'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.03'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.41'
$ws.Range('E28').Value = '  -1.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.273'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.44'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('E31').Value = '  -5.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.802.45'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.735'
$ws.Range('E33').Value = '  -2.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.87'
$ws.Range('E34').Value = '  +3.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9437'
$ws.Range('E35').Value = '  -4.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02781'
$ws.Range('E36').Value = '  -1.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2510'
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08810'
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.020'
$ws.Range('E39').Value = '  -4.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07080'
$ws.Range('E40').Value = '  -5.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.357'
$ws.Range('E41').Value = '  -3.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6986'
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.08'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.18'
$ws.Range('E44').Value = '  -4.97%  '
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6393'
$ws.Range('E46').Value = '  -3.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.299'
$ws.Range('E47').Value = '  -3.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.983'
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07956'
$ws.Range('E49').Value = '  -1.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.191'
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.88'
$ws.Range('E51').Value = '  -5.90%  '
